$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = '48.431.45'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(2, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.48%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.513.48'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(3, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.24%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.998'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(4, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.10%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '320.82'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(5, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.31%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '108.01'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(6, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.21%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(7, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.97%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(8, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.02%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -3.09%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '39.30'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(10, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.70%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '20.13'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(11, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +3.20%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0812'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(12, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.24%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(13, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.28%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.13'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(14, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.69%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.905.35'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(15, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.25%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.512.44'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(16, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.20%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.838'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(17, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.51%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '48.260.82'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(18, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.44%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(19, 2)
$cell.NumberFormat = "@"
$cell.Value = 'ImmutableX'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(19, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.02'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(19, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +9.78%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(20, 2)
$cell.NumberFormat = "@"
$cell.Value = 'InternetComputer(DFINITY)'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(20, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '13.12'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(20, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.52%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.70'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(21, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.35%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0944'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(22, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.26%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '71.66'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(23, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.01%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '274.56'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(24, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +11.00%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.55'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(25, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.59%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(26, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.07%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '26.05'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(27, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.03%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(28, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +2.14%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(29, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +5.59%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '9.82'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(30, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -3.61%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '35.35'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(31, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.46%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '49.69'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(32, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.20%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '19.40'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(33, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -3.39%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.34'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(34, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.92%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(35, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.18%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(36, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.12%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.97'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(37, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.23%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(38, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.82%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(39, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.45%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.111'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(40, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.66%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(41, 2)
$cell.NumberFormat = "@"
$cell.Value = 'WEMIXToken'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.22'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.72%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(42, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Monero'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '120.24'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.99%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '22.11'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(43, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.16%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0306'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(44, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +2.85%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.26'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(45, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +4.53%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.011.69'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(46, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.33%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.91'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(47, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +6.08%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.01'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(48, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.93%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '9.03'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(49, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.74%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.30'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(50, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +2.82%  '
$cell.Style = "Normal"

$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '79.80'
$cell.Style = "Normal"
$cell = $ws.Cells.Item(51, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +3.31%  '
$cell.Style = "Normal"

